$wb = $excel.ActiveWorkbook

# Rename the "aquisicoes" worksheet to "Folha1"
$ws = $wb.Worksheets.Item("aquisicoes")
$ws.Name = "Folha1"

# Update the selection/active cell on that sheet to F18
$ws.Activate()
$ws.Range("F18").Select()
